# 📊 Actualización automática del dashboard
# - Column J (comment_id): stored as inline text in the source export; re-typed as a
#   genuine number so it matches the other numeric ID columns (K, etc.).
# - Column M (fecha_comentario): number format switched from date-only to the same
#   date-time format already used by column L (created_time_processed).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Comentarios")

# comment_id values (row -> text value scraped from the API) to convert to numbers
$commentIds = [ordered]@{
    2  = "7603543574210986772"
    3  = "7603543320425923349"
    4  = "7603542496237093653"
    5  = "7603484564858192660"
    6  = "7603484464274785045"
    7  = "7603484333210682132"
    8  = "7602279642737066760"
    9  = "7602263847853064980"
    10  = "7602263706861667093"
    11  = "7601678572046107413"
    12  = "7601667808368165639"
    13  = "7601648158657413895"
    14  = "7601615764496007937"
    15  = "7601615659219469072"
    16  = "7601591880967619344"
    17  = "7601340708133749522"
    18  = "7601336152705352456"
    19  = "7601332782338245394"
    20  = "7601304306054611732"
    21  = "7601288389082530580"
    22  = "7601256810984932113"
    23  = "7601211223090266898"
    24  = "7601208873487926036"
    25  = "7601181962170630930"
    26  = "7601148384763151122"
    27  = "7601141061986370325"
    28  = "7601130474090464007"
    29  = "7601100678766134023"
    30  = "7601033369098896136"
    31  = "7600981628773204757"
    32  = "7600979599431025429"
    33  = "7600979263550817045"
    34  = "7600977406539612948"
    35  = "7600963433807495956"
    36  = "7600962044775301908"
    37  = "7600957864790197000"
    38  = "7600955753084863239"
    39  = "7600955758550123282"
    40  = "7600955752375370504"
    41  = "7600956815049229076"
    42  = "7600956719275000592"
}

foreach ($row in $commentIds.Keys) {
    $ws.Cells.Item($row, 10).Value = [double]$commentIds[$row]   # column J = comment_id
}

# fecha_comentario (column M) picks up the date-time number format used by column L
# (created_time_processed), replacing the old date-only format.
$lastRow = ($commentIds.Keys | Measure-Object -Maximum).Maximum
$firstRow = ($commentIds.Keys | Measure-Object -Minimum).Minimum
$ws.Range("M" + $firstRow + ":M" + $lastRow).NumberFormat = $ws.Range("L" + $firstRow).NumberFormat

